$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4561.905
$ws.Range("I116").Value = 3960
$ws.Range("J116").Value = 4932.3076
$ws.Range("K116").Value = 3960
$ws.Range("L116").Value = 4932.3076
$ws.Range("M116").Value = -518
$ws.Range("N116").Value = -11816.3076
$ws.Range("H132").Value = 2668688.5
$ws.Range("I132").Value = 2942690.2
$ws.Range("K132").Value = 8828070.600000001
$ws.Range("M132").Value = -8825540.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20688.367
$ws.Range("I32").Value = 17943.791
$ws.Range("J32").Value = 31666.666
$ws.Range("K32").Value = 17943.791
$ws.Range("L32").Value = 31666.666
$ws.Range("M32").Value = -17656.791
$ws.Range("N32").Value = -32240.666
$ws.Range("H44").Value = 22000
$ws.Range("J44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("N44").Value = -22976
$ws.Range("H55").Value = 26053
$ws.Range("J55").Value = 26053
$ws.Range("L55").Value = 26053
$ws.Range("N55").Value = -26683
$ws.Range("H61").Value = 3215.9
$ws.Range("I61").Value = 2557.4707
$ws.Range("J61").Value = 4076.923
$ws.Range("K61").Value = 2557.4707
$ws.Range("L61").Value = 4076.923
$ws.Range("M61").Value = -2345.4707
$ws.Range("N61").Value = -4500.923
$ws.Range("H97").Value = 725
$ws.Range("I97").Value = 610
$ws.Range("J97").Value = 1990
$ws.Range("K97").Value = 610
$ws.Range("L97").Value = 1990
$ws.Range("M97").Value = -114
$ws.Range("N97").Value = -2982
$ws.Range("H122").Value = 3900.25
$ws.Range("I122").Value = 2966.3333
$ws.Range("J122").Value = 4211.5557
$ws.Range("K122").Value = 8898.999899999999
$ws.Range("L122").Value = 12634.6671
$ws.Range("M122").Value = -6448.999899999999
$ws.Range("N122").Value = -17534.6671
$ws.Range("H132").Value = 2060.1428
$ws.Range("I132").Value = 1429.2572
$ws.Range("J132").Value = 3637.3572
$ws.Range("K132").Value = 4287.7716
$ws.Range("L132").Value = 10912.0716
$ws.Range("M132").Value = -1757.7716
$ws.Range("N132").Value = -15972.0716
$ws.Range("H133").Value = 28500
$ws.Range("J133").Value = 28500
$ws.Range("L133").Value = 28500
$ws.Range("N133").Value = -33560
$ws.Range("H136").Value = 3215.9
$ws.Range("I136").Value = 2557.4707
$ws.Range("J136").Value = 4076.923
$ws.Range("K136").Value = 7672.4121
$ws.Range("L136").Value = 12230.769
$ws.Range("M136").Value = -5122.4121
$ws.Range("N136").Value = -17330.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4089.9
$ws.Range("I99").Value = 3051.8
$ws.Range("J99").Value = 7204.2
$ws.Range("K99").Value = 3051.8
$ws.Range("L99").Value = 7204.2
$ws.Range("M99").Value = -1553.8
$ws.Range("N99").Value = -10200.2
$ws.Range("H122").Value = 30000
$ws.Range("J122").Value = 30000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -39800
$ws.Range("H132").Value = 29597.777
$ws.Range("J132").Value = 29597.777
$ws.Range("L132").Value = 29597.777
$ws.Range("N132").Value = -39717.777
$ws.Range("H134").Value = 3026.8
$ws.Range("I134").Value = 2823.2703
$ws.Range("K134").Value = 8469.8109
$ws.Range("M134").Value = -5934.8109

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1270.6428
$ws.Range("I22").Value = 558.9
$ws.Range("J22").Value = 3050
$ws.Range("K22").Value = 558.9
$ws.Range("L22").Value = 3050
$ws.Range("M22").Value = -208.9
$ws.Range("N22").Value = -3750
$ws.Range("H99").Value = 6502.3335
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -3502
$ws.Range("H105").Value = 2900.75
$ws.Range("I105").Value = 2556.4443
$ws.Range("J105").Value = 5999.5
$ws.Range("K105").Value = 2556.4443
$ws.Range("L105").Value = 5999.5
$ws.Range("M105").Value = -809.4443000000001
$ws.Range("N105").Value = -9493.5
$ws.Range("H126").Value = 6502.3335
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("H134").Value = 2689.6216
$ws.Range("I134").Value = 1918.4286
$ws.Range("J134").Value = 5088.8887
$ws.Range("K134").Value = 5755.2858
$ws.Range("L134").Value = 15266.6661
$ws.Range("M134").Value = -3220.2858
$ws.Range("N134").Value = -20336.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 9867.362999999999
$ws.Range("I6").Value = 434.42856
$ws.Range("J6").Value = 26375
$ws.Range("K6").Value = 1303.28568
$ws.Range("L6").Value = 79125
$ws.Range("M6").Value = -1190.28568
$ws.Range("N6").Value = -79351
$ws.Range("H11").Value = 17260.334
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 20612.4
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 61837.2
$ws.Range("M11").Value = -1360
$ws.Range("N11").Value = -62117.2
$ws.Range("H119").Value = 2897.2666
$ws.Range("I119").Value = 1243.1666
$ws.Range("K119").Value = 3729.4998
$ws.Range("M119").Value = 1108.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 30002000
$ws.Range("I3").Value = 150000000
$ws.Range("J3").Value = 2501
$ws.Range("K3").Value = 150000000
$ws.Range("L3").Value = 2501
$ws.Range("M3").Value = -149999884
$ws.Range("N3").Value = -2733
$ws.Range("H7").Value = 2666104.5
$ws.Range("I7").Value = 3153877
$ws.Range("J7").Value = 2032000.2
$ws.Range("K7").Value = 3153877
$ws.Range("L7").Value = 2032000.2
$ws.Range("M7").Value = -3153765
$ws.Range("N7").Value = -2032224.2
$ws.Range("H8").Value = 2666104.5
$ws.Range("I8").Value = 3153877
$ws.Range("J8").Value = 2032000.2
$ws.Range("K8").Value = 3153877
$ws.Range("L8").Value = 2032000.2
$ws.Range("M8").Value = -3153738
$ws.Range("N8").Value = -2032278.2
$ws.Range("H122").Value = 6806.75
$ws.Range("I122").Value = 1100
$ws.Range("J122").Value = 8709
$ws.Range("K122").Value = 3300
$ws.Range("L122").Value = 26127
$ws.Range("M122").Value = -850
$ws.Range("N122").Value = -31027
$ws.Range("H132").Value = 4802.032
$ws.Range("I132").Value = 5158.8096
$ws.Range("J132").Value = 4052.8
$ws.Range("K132").Value = 15476.4288
$ws.Range("L132").Value = 12158.4
$ws.Range("M132").Value = -12946.4288
$ws.Range("N132").Value = -17218.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1019.8421
$ws.Range("I22").Value = 380
$ws.Range("J22").Value = 1393.0834
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 1393.0834
$ws.Range("M22").Value = -85
$ws.Range("N22").Value = -1983.0834
$ws.Range("H27").Value = 1019.8421
$ws.Range("I27").Value = 380
$ws.Range("J27").Value = 1393.0834
$ws.Range("K27").Value = 380
$ws.Range("L27").Value = 1393.0834
$ws.Range("M27").Value = -273
$ws.Range("N27").Value = -1607.0834
$ws.Range("H122").Value = 4857.143
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("H123").Value = 29476.334
$ws.Range("J123").Value = 29476.334
$ws.Range("L123").Value = 29476.334
$ws.Range("N123").Value = -39276.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 60010
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 60010
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 60010
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -60470
$ws.Range("H122").Value = 1517.0834
$ws.Range("I122").Value = 1070.55
$ws.Range("J122").Value = 3749.75
$ws.Range("K122").Value = 3211.65
$ws.Range("L122").Value = 11249.25
$ws.Range("M122").Value = -761.6499999999996
$ws.Range("N122").Value = -16149.25
$ws.Range("H132").Value = 2515171
$ws.Range("I132").Value = 3231310.5
$ws.Range("K132").Value = 9693931.5
$ws.Range("M132").Value = -9691401.5
